$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (rows 2-7), columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T
# (columns F and L — detection rates — are unchanged)

$data = @{
    2 = @{ E=3; G=4.278141666666667; H=12.834425; I=0.9663225094340192; J=0.9663225094340191; K=3; M=0.5292956666666667; N=1.587887; O=0.135651968140022; P=0.1356519681400219; Q=2.264401845552778; R=20.379616609975; S=0.1310835502627296; T=0.1310835502627296 }
    3 = @{ E=3; G=4.278141666666667; H=12.834425; I=0.9663225094340192; J=0.9663225094340191; K=3; M=1.362890666666667; N=4.088672000000001; O=0.3492921120199358; P=0.3492921120199358; Q=5.830639348177779; R=52.47575413360001; S=0.3375288302126129; T=0.3375288302126129 }
    4 = @{ E=3; G=4.278141666666667; H=12.834425; I=0.9663225094340192; J=0.9663225094340191; K=3; M=2.009678666666666; N=6.029036; O=0.5150559198400423; P=0.5150559198400423; Q=8.597690040477776; R=77.37921036429999; S=0.4977101289586767; T=0.4977101289586767 }
    5 = @{ E=3; G=0.1490983333333333; H=0.447295; I=0.03367749056598092; J=0.03367749056598091; K=3; M=0.5292956666666667; N=1.587887; O=0.135651968140022; P=0.1356519681400219; Q=0.07891710174055556; R=0.7102539156650001; S=0.004568417877292333; T=0.004568417877292332 }
    6 = @{ E=3; G=0.1490983333333333; H=0.447295; I=0.03367749056598092; J=0.03367749056598091; K=3; M=1.362890666666667; N=4.088672000000001; O=0.3492921120199358; P=0.3492921120199358; Q=0.2032047269155556; R=1.82884254224; S=0.01176328180732294; T=0.01176328180732294 }
    7 = @{ E=3; G=0.1490983333333333; H=0.447295; I=0.03367749056598092; J=0.03367749056598091; K=3; M=2.009678666666666; N=6.029036; O=0.5150559198400423; P=0.5150559198400423; Q=0.2996397397355555; R=2.69675765762; S=0.01734579088136565; T=0.01734579088136564 }
}

foreach ($rowNum in $data.Keys) {
    $row = $data[$rowNum]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$rowNum").Value = $row[$col]
    }
}
